$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the added columns
$ws.Range("E1").Value = "eps_sigma_analytical"
$ws.Range("F1").Value = "eps_sigma_absolute"

# Fill constant values for all 27 data rows (rows 2-28)
$ws.Range("E2:E28").Value = 0.084195259341794
$ws.Range("F2:F28").Value = 0.3
